$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.777.16"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "1.865.50"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  +3.26%  "
$ws.Range("D5").Value = "324.50"
$ws.Range("E5").Value = "  +4.00%  "
$ws.Range("D6").Value = "1.036"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("D7").Value = "0.4426"
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("D9").Value = "0.07474"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("D11").Value = "21.74"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").Value = "1.882.22"
$ws.Range("E12").Value = "  -13.23%  "
$ws.Range("D13").Value = "5.565"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "6.771"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").Value = "0.07248"
$ws.Range("E15").Value = "  +4.41%  "
$ws.Range("D16").Value = "83.85"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "1.041"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "0.000009163"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").Value = "1.036"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").Value = "27.791.48"
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").Value = "5.320"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").Value = "11.38"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").Value = "1.997"
$ws.Range("E24").Value = "  +6.23%  "
$ws.Range("D25").Value = "158.90"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").Value = "18.88"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").Value = "5.335"
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("D28").Value = "1.986"
$ws.Range("E28").Value = "  +4.57%  "
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("D30").Value = "0.09070"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("D31").Value = "0.7788"
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("D32").Value = "3.106"
$ws.Range("E32").Value = "  +10.45%  "
$ws.Range("D33").Value = "1.213"
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("D35").Value = "1.037"
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("D36").Value = "1.156"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").Value = "0.01996"
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("D38").Value = "0.05356"
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("D39").Value = "2.874"
$ws.Range("E39").Value = "  +4.74%  "
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").Value = "0.1696"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("D42").Value = "6.895"
$ws.Range("E42").Value = "  +6.57%  "
$ws.Range("D43").Value = "8.692"
$ws.Range("E43").Value = "  +4.24%  "
$ws.Range("D44").Value = "109.93"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").Value = "10.72"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("E46").Value = "  +4.84%  "
$ws.Range("D47").Value = "0.4710"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").Value = "0.06476"
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("E49").Value = "  +4.24%  "
$ws.Range("D50").Value = "39.94"
$ws.Range("E50").Value = "  +3.64%  "
$ws.Range("D51").Value = "64.63"
$ws.Range("E51").Value = "  +2.76%  "